$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header Q8 in column J - copy formatting from the previous header (I1)
$ws.Range("I1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "Q8"

# Dates in column A stay the same; only the data in B:J changes.
# Row data: row -> array of values starting at column B
$rowData = @{
    2  = @(-4.324682701351703, 3.660091258637736, 1.731369225691766, 3.149790629511104, -3.635177286302406, 0.06544417180023943, -0.1340705272443026)
    3  = @(3.605598293642375, 1.676876260696405, 3.095297664515743, -3.689670251297767, 0.01095120680487852, -0.1885634922396635)
    4  = @(0.8917107665847916, 2.310132170404129, -4.47483574540938, -0.7742142873067352, -0.9737289863512772, 1.180395370885164, -3.253940853883919, 1.466850648880309, -0.1948090203170254)
    5  = @(2.584158590043411, -4.200809325770098, -0.5001878676674536, -0.6997025667119956, 1.454421790524445, -2.979914434244638, 1.74087706851959, 0.07921739932225619)
    6  = @(-4.14234893673768, -0.4417274786350349, -0.6412421776795769, 1.512882179556864, -2.921454045212219, 1.799337457552009, 0.1376777883546749)
    7  = @(-0.6288648973096044, -0.8283795963541464, 1.325744760882294, -3.108591463886789, 1.61220003887744, -0.04945963031989459)
    8  = @(-0.9420242164422319, 1.212100140794209, -3.222236083974874, 1.498555418789354, -0.1631042504079802, -1.553285820388052, -1.037158441845645, 1.563884698290838)
    9  = @(0.9834580318344777, -3.450878192934606, 1.269913309829623, -0.3917463593677113, -1.781927929347783, -1.265800550805376, 1.335242589331107)
    10 = @(-3.353876282308941, 1.366915220455287, -0.294744448742047, -1.684926018722119, -1.168798640179712, 1.432244499956771)
    11 = @(1.609820403724652, -0.05183926547268242, -1.442020835452754, -0.9258934569103474, 1.675149683226136)
    12 = @(-0.7294052785381211, -2.119586848518193, -1.603459469975786, 0.997583670160697)
    13 = @(-1.846431246597561, -1.330303868055154, 1.270739272081329)
    14 = @(-1.121732059169287, 1.479311080967196)
    15 = @(1.884604928957667)
    16 = @()
}

foreach ($r in 2..16) {
    $vals = $rowData[$r]
    # Clear any previous values in B:J for this row first
    $ws.Range("B$r`:J$r").ClearContents()
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
